$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from column A into the new H:I columns for every
# used row, matching header (row 1) vs. body row (rows 2-22) styling.
$ws.Range("A1:A22").Copy($ws.Range("H1:H22"))
$ws.Range("A1:A22").Copy($ws.Range("I1:I22"))

# Header row
$ws.Range("H1").Value = "EUC"
$ws.Range("I1").Value = "UIP"

# EUC column values per row (UIP column stays blank)
$ws.Range("H2").Value = "Done"
$ws.Range("H3").Value = "Done"
$ws.Range("H4").Value = "Done"
$ws.Range("H5").Value = "Done"
$ws.Range("H6").Value = "Done"
$ws.Range("H7").Value = "Done"
$ws.Range("H8").Value = "Done"
$ws.Range("H9").Value = "Not Started"
$ws.Range("H14").Value = "Not Started"
$ws.Range("H15").Value = "Not Started"
$ws.Range("H16").Value = "Not Started"
$ws.Range("H20").Value = "Not Started"
$ws.Range("H21").Value = "Not Started"
$ws.Range("H22").Value = "Not Started"

# I column (UIP) has no values beyond the header - clear any copied content
$ws.Range("I2:I22").ClearContents()

# Rows with no EUC status stay blank - clear the values copied from column A
$ws.Range("H10:H13").ClearContents()
$ws.Range("H17:H19").ClearContents()

# Update the current selection to match the authored file (J1)
$ws.Range("J1").Select()
